$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Insert two new rows before row 52 to grow the table from 52 to 54 data rows
$ws1.Rows.Item(52).Insert()
$ws1.Rows.Item(52).Insert()

# --- Sheet "Recommandations": rewrite rows 2-54 (header row 1 is unchanged) ---
$recoData = @{}
$recoData[2] = @('BRVM - SERVICES PUBLICS', 0, 8, 3359.91, 107.92, '🟡 Observer', '➖ Neutre')
$recoData[3] = @('AIR LIQUIDE CI', 0, 4, 2765, 685, '🟡 Observer', '➖ Neutre')
$recoData[4] = @('NEI-CEDA CI', 0, 4, 2715, 685, '🟡 Observer', '➖ Neutre')
$recoData[5] = @('BRVM - AUTRES SECTEURS', 0, 4, 2587.1, 641.23, '🟡 Observer', '➖ Neutre')
$recoData[6] = @('BRVM - DISTRIBUTION', 0, 4, 1677.06, 432.61, '🟡 Observer', '➖ Neutre')
$recoData[7] = @('BRVM - TRANSPORT', 0, 4, 1468.64, 367.16, '🟡 Observer', '➖ Neutre')
$recoData[8] = @('BRVM - AGRICULTURE', 0, 4, 1362.62, 340.44, '🟡 Observer', '➖ Neutre')
$recoData[9] = @('CFAO MOTORS CI', 0, 1, 940, 940, '🟡 Observer', '➖ Neutre')
$recoData[10] = @('SETAO CI', 0, 1, 935, 935, '🟡 Observer', '➖ Neutre')
$recoData[11] = @('BRVM - INDUSTRIELS', 0, 4, 566.24, 137.74, '🟡 Observer', '➖ Neutre')
$recoData[12] = @('BRVM-PRESTIGE', 0, 4, 547.71, 137.06, '🟡 Observer', '➖ Neutre')
$recoData[13] = @('BRVM - CONSOMMATION DISCRETIONNAIRE', 0, 4, 546.09, 143.47, '🟡 Observer', '➖ Neutre')
$recoData[14] = @('BRVM - FINANCES', 0, 4, 537.08, 134.51, '🟡 Observer', '➖ Neutre')
$recoData[15] = @('BRVM - SERVICES FINANCIERS', 0, 4, 527.8200000000001, 132.19, '🟡 Observer', '➖ Neutre')
$recoData[16] = @('BRVM - ENERGIE', 0, 4, 430.44, 107.36, '🟡 Observer', '➖ Neutre')
$recoData[17] = @('BRVM - TELECOMMUNICATIONS', 0, 4, 381.3, 95.38, '🟡 Observer', '➖ Neutre')
$recoData[18] = @('BRVM - CONSOMMATION DE BASE            (**)', 0, 2, 371.92, 185.85, '🟡 Observer', '➖ Neutre')
$recoData[19] = @('BRVM - INDUSTRIE                 (**)', 0, 1, 216.31, 216.31, '🟡 Observer', '➖ Neutre')
$recoData[20] = @('BRVM - INDUSTRIE                       (**)', 0, 1, 214.57, 214.57, '🟡 Observer', '➖ Neutre')
$recoData[21] = @('BRVM - INDUSTRIE                         (**)', 0, 1, 214.08, 214.08, '🟡 Observer', '➖ Neutre')
$recoData[22] = @('BRVM - INDUSTRIE                      (**)', 0, 1, 211.51, 211.51, '🟡 Observer', '➖ Neutre')
$recoData[23] = @('BRVM-PRINCIPAL                         (**)', 0, 1, 192.48, 192.48, '🟡 Observer', '➖ Neutre')
$recoData[24] = @('BRVM-PRINCIPAL                           (**)', 0, 1, 192.4, 192.4, '🟡 Observer', '➖ Neutre')
$recoData[25] = @('BRVM-PRINCIPAL                          (**)', 0, 1, 192.26, 192.26, '🟡 Observer', '➖ Neutre')
$recoData[26] = @('BRVM-PRINCIPAL                    (**)', 0, 1, 192.05, 192.05, '🟡 Observer', '➖ Neutre')
$recoData[27] = @('BRVM - CONSOMMATION DE BASE        (**)', 0, 1, 184.97, 184.97, '🟡 Observer', '➖ Neutre')
$recoData[28] = @('BRVM - CONSOMMATION DE BASE             (**)', 0, 1, 183.77, 183.77, '🟡 Observer', '➖ Neutre')
$recoData[29] = @('UNIWAX CI (UNXC)', 4, 0, 25.75, 4.62, '🟢 Achat', '✅ Renforcer')
$recoData[30] = @('ORAGROUP TOGO (ORGT)', 3, 0, 22.38, 7.43, '🟢 Achat', '✅ Renforcer')
$recoData[31] = @('CFAO MOTORS CI (CFAC)', 3, 0, 18.72, 5.53, '🟢 Achat', '✅ Renforcer')
$recoData[32] = @('BERNABE CI (BNBC)', 1, 0, 7.07, 7.07, '🟡 Observer', '➖ Neutre')
$recoData[33] = @('SAFCA CI (SAFC)', 1, 0, 6.49, 6.49, '🟡 Observer', '➖ Neutre')
$recoData[34] = @('SETAO CI (STAC)', 2, 1, 5.32, -7.44, '🟡 Observer', '👀 À surveiller')
$recoData[35] = @('CORIS BANK INTERNATIONAL (CBIBF)', 1, 0, 4.52, 4.52, '🟡 Observer', '➖ Neutre')
$recoData[36] = @('NSIA BANQUE COTE D''IVOIRE (NSBC)', 1, 1, 4.5, -2.96, '🟡 Observer', '👀 À surveiller')
$recoData[37] = @('BANK OF AFRICA BN (BOAB)', 1, 0, 3.7, 3.7, '🟡 Observer', '➖ Neutre')
$recoData[38] = @('BANK OF AFRICA ML (BOAM)', 1, 0, 3.49, 3.49, '🟡 Observer', '➖ Neutre')
$recoData[39] = @('BICI CI (BICC)', 1, 0, 1.96, 1.96, '🟡 Observer', '➖ Neutre')
$recoData[40] = @('SAPH CI (SPHC)', 1, 1, 1.47, -1.72, '🟡 Observer', '👀 À surveiller')
$recoData[41] = @('TOTAL', 0, 4, 0, 0, '🟡 Observer', '➖ Neutre')
$recoData[42] = @('PALM CI (PALC)', 0, 1, -1.1, -1.1, '🟡 Observer', '➖ Neutre')
$recoData[43] = @('NESTLE CI (NTLC)', 0, 1, -1.15, -1.15, '🟡 Observer', '➖ Neutre')
$recoData[44] = @('AIR LIQUIDE CI (SIVC)', 0, 1, -1.44, -1.44, '🟡 Observer', '➖ Neutre')
$recoData[45] = @('ORANGE COTE D''IVOIRE (ORAC)', 0, 1, -2.03, -2.03, '🟡 Observer', '➖ Neutre')
$recoData[46] = @('SOLIBRA CI (SLBC)', 0, 1, -3.1, -3.1, '🟡 Observer', '➖ Neutre')
$recoData[47] = @('SOCIETE IVOIRIENNE DE BANQUE  (SIBC)', 0, 1, -3.42, -3.42, '🟡 Observer', '➖ Neutre')
$recoData[48] = @('BANK OF AFRICA SENEGAL (BOAS)', 0, 1, -3.7, -3.7, '🟡 Observer', '➖ Neutre')
$recoData[49] = @('SICOR CI (SICC)', 0, 1, -4, -4, '🟡 Observer', '➖ Neutre')
$recoData[50] = @('SUCRIVOIRE (SCRC)', 0, 2, -4.42, -2.6, '🟡 Observer', '➖ Neutre')
$recoData[51] = @('NEI-CEDA CI (NEIC)', 0, 1, -5.84, -5.84, '🟡 Observer', '➖ Neutre')
$recoData[52] = @('FILTISAC CI (FTSC)', 0, 2, -12.19, -7.41, '🟡 Observer', '➖ Neutre')
$recoData[53] = @('BANK OF AFRICA BF (BOABF)', 0, 2, -12.33, -4.91, '🟡 Observer', '➖ Neutre')
$recoData[54] = @('UNILEVER CI (UNLC)', 0, 2, -14.97, -7.49, '🟡 Observer', '➖ Neutre')

foreach ($r in $recoData.Keys) {
    $vals = $recoData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws1.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}

# --- Sheet "Top_YTD": update labels/values for rows 2-11 (header row 1 unchanged) ---
$ytdData = @{}
$ytdData[2] = @('BRVM - SERVICES PUBLICS', 9027952.039999999)
$ytdData[3] = @('AIR LIQUIDE CI', 391850.3)
$ytdData[4] = @('NEI-CEDA CI', 367447.4)
$ytdData[5] = @('BRVM - AUTRES SECTEURS', 310887.06)
$ytdData[6] = @('BRVM - DISTRIBUTION', 72543.82000000001)
$ytdData[7] = @('BRVM - TRANSPORT', 47526.72)
$ytdData[8] = @('BRVM - AGRICULTURE', 37602.14)
$ytdData[9] = @('BRVM - INDUSTRIELS', 3304.27)
$ytdData[10] = @('BRVM-PRESTIGE', 3051.08)
$ytdData[11] = @('BRVM - CONSOMMATION DISCRETIONNAIRE', 3026.56)

foreach ($r in $ytdData.Keys) {
    $vals = $ytdData[$r]
    $ws2.Cells.Item($r, 1).Value = $vals[0]
    $ws2.Cells.Item($r, 2).Value = $vals[1]
}
